# Auto-generated script to reorder gene-set string literals in column I (Intersection)
# across all 6 cluster worksheets, per the target diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("enrichment_cluster_1")
$ws.Range("I2").Value = '{''ERCC1'', ''CDK7'', ''TERT'', ''RAD52'', ''XRCC5'', ''LMNA'', ''RAD51'', ''AKT1'', ''PTEN'', ''ATR'', ''POLA1'', ''ERCC2'', ''CDKN1A'', ''ATM'', ''WRN'', ''TP53'', ''ERCC8'', ''BRCA1''}'
$ws.Range("I5").Value = '{''APOO'', ''AR'', ''TERF1'', ''RAD51D'', ''ZMPSTE24''}'
$ws.Range("I6").Value = '{''TDP1'', ''RAD9A'', ''NEIL1'', ''LIG1'', ''POLE'', ''DVL1'', ''DVL2'', ''DAXX'', ''MCM7''}'
$ws.Range("I7").Value = '{''DCLRE1C'', ''RBBP8'', ''CDKN1A'', ''HADHB'', ''CDK2''}'
$ws.Range("I8").Value = '{''MCM4'', ''RAD51C'', ''RAD50'', ''TP53BP1'', ''DCP2'', ''CCND2'', ''USP12'', ''CTNNB1'', ''ATRX''}'
$ws.Range("I9").Value = '{''XRCC6'', ''CUL4A'', ''SUMO2'', ''CDKN1A'', ''USP1''}'
$ws.Range("I10").Value = '{''MTX2'', ''AR'', ''SAMM50'', ''CDK2'', ''CHCHD10''}'
$ws.Range("I11").Value = '{''MDM4'', ''USP34'', ''CCNH''}'
$ws.Range("I12").Value = '{''RAD9A'', ''POLD1'', ''KDM1A'', ''CTNNB1'', ''DVL2'', ''TP53'', ''RTEL1''}'
$ws = $wb.Worksheets.Item("enrichment_cluster_3")
$ws.Range("I2").Value = '{''RB1'', ''CDK7'', ''AKT1'', ''PML'', ''SIRT1'', ''TFDP1'', ''FOXO1'', ''FOXO3'', ''E2F1'', ''ERCC2'', ''CDKN1A'', ''CDKN2B'', ''TP53'', ''CDKN2A'', ''BRCA1'', ''NCOR2''}'
$ws.Range("I5").Value = '{''FOXO1'', ''FOXO3'', ''TEAD1'', ''TBX5'', ''MAML2'', ''ERLEC1'', ''WWTR1'', ''AMOTL2'', ''ELK3'', ''RBPJ'', ''APH1B'', ''ARHGAP26'', ''CDH10'', ''ZFPM2'', ''UBE2J1'', ''AMOT'', ''THRB'', ''ESR1'', ''AR''}'
$ws.Range("I6").Value = '{''GTF3C1'', ''BRF1'', ''CDKN2D'', ''KAT2A'', ''VDR'', ''RXRB'', ''EGR2'', ''GTF2H4'', ''DAXX'', ''HDAC10''}'
$ws.Range("I7").Value = '{''MAML2'', ''NOTCH2'', ''YAP1'', ''EGFR'', ''SOX9'', ''PKN3'', ''SIRT1'', ''CBFB'', ''STAT3'', ''CDKN1A'', ''SMAD1'', ''CDK2'', ''BOC'', ''CDKN2A''}'
$ws.Range("I8").Value = '{''GTF2H1'', ''SMAD3'', ''THRB'', ''SMARCA2'', ''RXRB'', ''CCND2'', ''CTNNB1'', ''RANGAP1'', ''SRC''}'
$ws.Range("I9").Value = '{''SMAD4'', ''RBPJ'', ''SNW1'', ''CREBBP'', ''ADAM10'', ''SIRT1'', ''FOXO1'', ''SUMO2'', ''PIK3R1'', ''FOXO3'', ''CDKN1A'', ''CTNNA1'', ''CCNT1''}'
$ws.Range("I10").Value = '{''AR'', ''PML'', ''CDK2'', ''CDKN2C'', ''HIVEP3''}'
$ws.Range("I11").Value = '{''CDKN2A'', ''MAF'', ''CCNH'', ''SEL1L''}'
$ws.Range("I12").Value = '{''SMARCA4'', ''ELK1'', ''MED12'', ''TGIF1'', ''SMARCB1'', ''SMARCC2'', ''KDM1A'', ''CTNNB1'', ''TP53'', ''JUP'', ''VGLL4''}'
$ws = $wb.Worksheets.Item("enrichment_cluster_5")
$ws.Range("I2").Value = '{''CDC42'', ''HRAS'', ''IGF1R'', ''ELN'', ''INSR'', ''HTRA2'', ''IRS1'', ''SHC1'', ''PIK3CB''}'
$ws.Range("I3").Value = '{''PAPSS2'', ''RET'', ''PDGFB'', ''SYTL3''}'
$ws.Range("I4").Value = '{''FBLN2'', ''CR2'', ''MS4A1'', ''PYCR1''}'
$ws.Range("I5").Value = '{''RCN1'', ''COL5A2'', ''AOX1'', ''LOX'', ''PIK3R3'', ''PTK2'', ''GPX7'', ''EDN3'', ''PDGFRA'', ''CYCS'', ''PDGFRB'', ''IRS1'', ''CERS6'', ''CYS1'', ''EMILIN1'', ''KIF5B'', ''THBS1'', ''FYN'', ''DLC1'', ''RPGRIP1L'', ''ITGA1'', ''C3'', ''UACA'', ''MFAP5'', ''ADAMTS2'', ''ABL2'', ''LTBP1'', ''SIRPA'', ''CXCL12'', ''JAM3'', ''HSP90AA1'', ''RHOA'', ''FBN1'', ''ELN'', ''RHOJ'', ''SH3KBP1'', ''ITGA9'', ''MSR1'', ''WIPF1'', ''MFAP4'', ''GPX8'', ''ESR1'', ''ANKS6'', ''NRP1'', ''GALNT12'', ''TGFBR1'', ''FBLN5'', ''DOCK4'', ''RHOQ'', ''ITGAV'', ''PIK3AP1''}'
$ws.Range("I6").Value = '{''ITGB7'', ''GRB7'', ''SEMA4D'', ''MOGS'', ''DVL2'', ''EPHB3'', ''PLK1'', ''COL2A1'', ''CNKSR1'', ''PLEKHH3'', ''GJB6'', ''PIK3R2'', ''PSAT1'', ''TROAP'', ''COCH'', ''PLXNB1'', ''HGS'', ''ADAM15''}'
$ws.Range("I7").Value = '{''HGF'', ''SERPINA1'', ''LCP2'', ''CD2AP'', ''P4HA1'', ''FGG'', ''EGFR'', ''FYN'', ''PKN3'', ''TGFB1'', ''C3'', ''RAB27A'', ''PAPSS2'', ''LTBP1'', ''ITSN1'', ''LRP4'', ''ITGB1'', ''RANBP3L'', ''FCGR1A'', ''CD44'', ''MSR1'', ''SHC1'', ''TLN1'', ''FLT1'', ''PARD3'', ''TGFB3'', ''CD81'', ''SYK'', ''TGFBR1'', ''SPP1'', ''DOCK4'', ''ITGB5'', ''PXN'', ''RHOQ'', ''ECE1'', ''ITGAV'', ''PIK3AP1''}'
$ws.Range("I8").Value = '{''SOS1'', ''MET'', ''DLG2'', ''SRC'', ''PIK3CB'', ''DOCK3'', ''F11R'', ''DLG1'', ''CTNNB1'', ''EPHA5'', ''YWHAB'', ''CHRNA7'', ''GOPC'', ''AGPAT1'', ''MYO5A'', ''RCN2'', ''AKAP10'', ''DLGAP1'', ''SMAD3''}'
$ws.Range("I9").Value = '{''SLC44A2'', ''YWHAZ'', ''LRCH3'', ''KRAS'', ''XIAP'', ''COL4A3'', ''SMARCA5'', ''COL4A5'', ''PIK3R1'', ''EIF4H'', ''YWHAB''}'
$ws.Range("I10").Value = '{''AGPAT1'', ''PARD6A'', ''PDGFB'', ''VEGFC'', ''ALKBH7'', ''EML4''}'
$ws.Range("I11").Value = '{''MFAP4'', ''MYO5A'', ''YWHAZ'', ''GPX8'', ''MITF'', ''KITLG'', ''DLC1'', ''EML1'', ''XIAP'', ''RHOQ'', ''OSTF1'', ''TFRC''}'
$ws.Range("I12").Value = '{''EPHB4'', ''MED12'', ''PARD6G'', ''UNC119B'', ''AGRN'', ''CTNNB1'', ''DVL2'', ''GTF2IRD1'', ''PLXNB1'', ''SHC1'', ''ITGA3''}'
$ws = $wb.Worksheets.Item("enrichment_cluster_4")
$ws.Range("I2").Value = '{''NFKBIA'', ''HTT'', ''HSPA1A'', ''AKT1'', ''FOXO3'', ''IKBKB'', ''TP53''}'
$ws.Range("I5").Value = '{''CD14'', ''CFLAR'', ''IRAK3'', ''HSP90AA1'', ''RAB8B'', ''CASP10'', ''BCL2'', ''FOXO3'', ''TLR4'', ''FKBP5''}'
$ws.Range("I6").Value = '{''NFKB2'', ''NLRP1'', ''TBC1D17'', ''IKBKE'', ''TAB1'', ''IKBKB'', ''LRSAM1''}'
$ws.Range("I7").Value = '{''CASP1'', ''NFKBIA'', ''CD14'', ''BCL2L11'', ''TNIP2'', ''TNFRSF1B'', ''RELA'', ''SHMT1'', ''SFN'', ''BCL2'', ''LTBR'', ''MCL1'', ''FKBP5'', ''TLR4'', ''TNFRSF1A'', ''MAP3K8'', ''MYD88''}'
$ws.Range("I9").Value = '{''CUL4A'', ''CREBBP'', ''FOXO3''}'
$ws.Range("I12").Value = '{''SHMT2'', ''TP53'', ''FADD'', ''TRAF2''}'
$ws = $wb.Worksheets.Item("enrichment_cluster_2")
$ws.Range("I2").Value = '{''RB1'', ''CDK7'', ''AKT1'', ''TFDP1'', ''FOXO1'', ''MYC'', ''FOXO3'', ''E2F1'', ''CDKN1A'', ''CDKN2B'', ''TP53'', ''ERCC8'', ''CDKN2A''}'
$ws.Range("I3").Value = '{''DDB2'', ''CDKN2A'', ''SYTL3''}'
$ws.Range("I5").Value = '{''NQO1'', ''AR'', ''ESR1'', ''KIF5B'', ''FOXO1'', ''FOXO3'', ''YAF2''}'
$ws.Range("I7").Value = '{''GMNN'', ''STAT3'', ''CDKN1A'', ''SMAD1'', ''RAB27A'', ''CDK2'', ''CDKN2A''}'
$ws.Range("I8").Value = '{''MYO5A'', ''MCM4'', ''OAZ2'', ''YY1'', ''CDC6'', ''KCTD1'', ''CCND2'', ''CTNNB1'', ''SMAD3'', ''YAF2''}'
$ws.Range("I9").Value = '{''SMAD4'', ''PHC3'', ''YY1'', ''CREBBP'', ''CUL4A'', ''FOXO1'', ''FOXO3'', ''EPC1'', ''CDKN1A'', ''EHMT1''}'
$ws.Range("I10").Value = '{''AR'', ''CDK2'', ''CDKN2C'', ''CITED4''}'
$ws.Range("I11").Value = '{''MYO5A'', ''CCNH'', ''CDKN2A''}'
$ws.Range("I12").Value = '{''EHMT2'', ''PHC1'', ''SMARCA4'', ''TGIF1'', ''KDM1A'', ''CTNNB1'', ''TP53'', ''CITED4''}'
$ws = $wb.Worksheets.Item("enrichment_cluster_6")
$ws.Range("I4").Value = '{''NT5E'', ''CHML'', ''PAWR''}'
$ws.Range("I5").Value = '{''RABEP1'', ''OSBPL1A'', ''RAB32'', ''RAB8B'', ''ABCD2'', ''APH1B'', ''SALL1'', ''PEX19'', ''PTP4A2'', ''ARHGDIG'', ''NT5E'', ''PRKRA'', ''LONRF3''}'
$ws.Range("I6").Value = '{''SLC6A8'', ''C19orf54'', ''TMEM102'', ''S1PR5'', ''PPFIA3'', ''TBC1D17'', ''SIRT2'', ''RAB38'', ''GGA1''}'
$ws.Range("I7").Value = '{''SLC6A8'', ''TRIM4'', ''RAB34'', ''ZFP36'', ''USP3'', ''RAB13'', ''PAWR'', ''RAB31''}'
$ws.Range("I8").Value = '{''RAB35'', ''CHM'', ''RAB14'', ''RABEP1'', ''OSBPL1A'', ''ABCD2'', ''VAMP2'', ''TMEM9'', ''BACE1'', ''PRKRA'', ''RAB1A''}'
$ws.Range("I9").Value = '{''DICER1'', ''CHM'', ''RAB10'', ''ZFP36'', ''PRKRA''}'
$ws.Range("I10").Value = '{''GGA3'', ''TMEM9'', ''C19orf54'', ''EHD3''}'
$ws.Range("I11").Value = '{''NT5E'', ''RAB3GAP1'', ''RAB31''}'
$ws.Range("I12").Value = '{''MICALL1'', ''TMEM102'', ''PPFIA3''}'
